$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column L: header "break_on_off" and mark "break" trials
$ws.Range("L1").Value = "break_on_off"
$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").HorizontalAlignment = -4108

$lastRow = 73
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 12).Value = 0
}

# Mark rows where a break occurred (per target data)
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(37, 12).Value = 1
$ws.Cells.Item(54, 12).Value = 1

# Update selection to match target view
$ws.Range("L1:L73").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
